# Re-sort the comma-separated "Recorded By" names in column G using a
# case-sensitive (ordinal/ASCII) sort order, e.g.
#   "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# because uppercase letters (e.g. 'S' = 83) sort before lowercase letters
# (e.g. 'a'.. = 97+), which matches how the source data was re-generated.

function Compare-Ordinal($left, $right) {
    $lenLeft = $left.Length
    $lenRight = $right.Length
    $minLen = [Math]::Min($lenLeft, $lenRight)
    for ($charIdx = 0; $charIdx -lt $minLen; $charIdx++) {
        $codeLeft = [int][char]$left.Substring($charIdx, 1)
        $codeRight = [int][char]$right.Substring($charIdx, 1)
        if ($codeLeft -lt $codeRight) { return -1 }
        if ($codeLeft -gt $codeRight) { return 1 }
    }
    if ($lenLeft -lt $lenRight) { return -1 }
    if ($lenLeft -gt $lenRight) { return 1 }
    return 0
}

function Sort-Ordinal($items) {
    $total = $items.Length
    for ($passIdx = 0; $passIdx -lt $total; $passIdx++) {
        for ($cmpIdx = 0; $cmpIdx -lt ($total - $passIdx - 1); $cmpIdx++) {
            $cmpResult = Compare-Ordinal $items[$cmpIdx] $items[$cmpIdx + 1]
            if ($cmpResult -gt 0) {
                $swapTmp = $items[$cmpIdx]
                $items[$cmpIdx] = $items[$cmpIdx + 1]
                $items[$cmpIdx + 1] = $swapTmp
            }
        }
    }
    return $items
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($rowNum = 2; $rowNum -le $lastRow; $rowNum++) {
    $cell = $ws.Cells.Item($rowNum, 7)
    $cellValue = $cell.Value2
    if ($cellValue -eq $null) { continue }
    if ($cellValue -eq "") { continue }

    $nameParts = $cellValue -split ", "
    if ($nameParts.Length -gt 1) {
        $sortedParts = Sort-Ordinal $nameParts
        $newValue = [string]::Join(", ", $sortedParts)
        if ($newValue -ne $cellValue) {
            $cell.Value2 = $newValue
        }
    }
}
